$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID19_TIMESERIESDATA")

# New header cell Y1 - date label "10-10-2020" stored as text (like X1),
# with the same border/bold/alignment formatting as the preceding date header cells.
$ws.Range("Y1").NumberFormat = "@"
$ws.Range("Y1").Value = "10-10-2020"
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)

# New data column Y values (row 2 .. row 36), one new deceased-count column
# added to the time series (date 10-10-2020).
$values = @(
    55,
    6159,
    23,
    802,
    934,
    188,
    1196,
    2,
    5692,
    491,
    3547,
    1562,
    245,
    1306,
    781,
    9789,
    955,
    63,
    2574,
    39732,
    86,
    61,
    0,
    17,
    991,
    558,
    3773,
    1621,
    53,
    10120,
    1217,
    313,
    716,
    6293,
    5501
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 25).Value = $values[$i]
}
